# Update the cryptos price/volume table with the latest scraped values.
# Column D = Price, Column E = Volume(1h). Row numbers below are the
# worksheet row numbers (header is row 1, data starts at row 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "33.951.99";  NeedsText = $false; E = "  +9.31%  " },
    @{ Row = 3;  D = "1.784.29";   NeedsText = $false; E = "  +5.98%  " },
    @{ Row = 4;  D = $null;        NeedsText = $false; E = "  +0.11%  " },
    @{ Row = 5;  D = "225.03";     NeedsText = $true;  E = "  +2.21%  " },
    @{ Row = 6;  D = $null;        NeedsText = $false; E = "  +4.46%  " },
    @{ Row = 7;  D = "0.999";      NeedsText = $true;  E = "  +0.05%  " },
    @{ Row = 8;  D = "30.75";      NeedsText = $true;  E = "  +5.22%  " },
    @{ Row = 9;  D = "46.28";      NeedsText = $true;  E = "  +4.05%  " },
    @{ Row = 10; D = $null;        NeedsText = $false; E = "  +4.78%  " },
    @{ Row = 11; D = "0.0662";     NeedsText = $true;  E = "  +3.31%  " },
    @{ Row = 12; D = "0.0924";     NeedsText = $true;  E = "  +1.85%  " },
    @{ Row = 13; D = "2.038.74";   NeedsText = $false; E = "  +5.80%  " },
    @{ Row = 14; D = "1.787.85";   NeedsText = $false; E = "  +6.28%  " },
    @{ Row = 15; D = "0.631";      NeedsText = $true;  E = "  +4.07%  " },
    @{ Row = 16; D = "33.922.87";  NeedsText = $false; E = "  +9.15%  " },
    @{ Row = 17; D = "10.02";      NeedsText = $true;  E = "  -2.10%  " },
    @{ Row = 18; D = "4.20";       NeedsText = $true;  E = "  +1.27%  " },
    @{ Row = 19; D = "68.70";      NeedsText = $true;  E = "  +3.09%  " },
    @{ Row = 20; D = "252.07";     NeedsText = $true;  E = "  +1.75%  " },
    @{ Row = 21; D = "0.0₃0741";   NeedsText = $false; E = "  +2.76%  " },
    @{ Row = 22; D = "1.00";       NeedsText = $true;  E = "  +0.17%  " },
    @{ Row = 23; D = "10.32";      NeedsText = $true;  E = "  +3.00%  " },
    @{ Row = 24; D = "4.22";       NeedsText = $true;  E = "  -1.55%  " },
    @{ Row = 25; D = $null;        NeedsText = $false; E = "  -0.20%  " },
    @{ Row = 26; D = "158.42";     NeedsText = $true;  E = "  -0.09%  " },
    @{ Row = 27; D = "16.49";      NeedsText = $true;  E = "  +3.54%  " },
    @{ Row = 28; D = $null;        NeedsText = $false; E = "  +1.95%  " },
    @{ Row = 29; D = "6.96";       NeedsText = $true;  E = "  +3.84%  " },
    @{ Row = 30; D = $null;        NeedsText = $false; E = "  +0.09%  " },
    @{ Row = 31; D = $null;        NeedsText = $false; E = "  +8.18%  " },
    @{ Row = 32; D = "0.0515";     NeedsText = $true;  E = "  +3.17%  " },
    @{ Row = 33; D = $null;        NeedsText = $false; E = "  +4.20%  " },
    @{ Row = 34; D = "3.56";       NeedsText = $true;  E = "  +6.51%  " },
    @{ Row = 35; D = "1.490.04";   NeedsText = $false; E = "  -1.61%  " },
    @{ Row = 36; D = "1.80";       NeedsText = $true;  E = "  +3.42%  " },
    @{ Row = 37; D = $null;        NeedsText = $false; E = "  +2.90%  " },
    @{ Row = 38; D = "0.633";      NeedsText = $true;  E = "  +3.60%  " },
    @{ Row = 39; D = $null;        NeedsText = $false; E = "  +2.92%  " },
    @{ Row = 40; D = "83.38";      NeedsText = $true;  E = "  -0.65%  " },
    @{ Row = 41; D = "2.35";       NeedsText = $true;  E = "  +2.59%  " },
    @{ Row = 42; D = "2.70";       NeedsText = $true;  E = "  +1.15%  " },
    @{ Row = 43; D = "0.888";      NeedsText = $true;  E = "  +5.74%  " },
    @{ Row = 44; D = "2.09";       NeedsText = $true;  E = "  +2.79%  " },
    @{ Row = 45; D = "0.0510";     NeedsText = $true;  E = "  +1.43%  " },
    @{ Row = 46; D = $null;        NeedsText = $false; E = "  +3.10%  " },
    @{ Row = 47; D = "1.936.43";   NeedsText = $false; E = "  +6.35%  " },
    @{ Row = 48; D = $null;        NeedsText = $false; E = "  +3.74%  " },
    @{ Row = 49; D = $null;        NeedsText = $false; E = "  +0.18%  " },
    @{ Row = 50; D = "11.90";      NeedsText = $true;  E = "  +14.49%  " },
    @{ Row = 51; D = "50.81";      NeedsText = $true;  E = "  -2.31%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($u.Row, 4)
        if ($u.NeedsText) {
            # These prices look like plain numbers ("225.03", "1.00", ...).
            # Force the cell to Text first so Excel doesn't silently convert
            # the string into a numeric value (which would drop formatting
            # like trailing zeros, e.g. "1.00" -> 1).
            $cell.NumberFormat = "@"
        }
        $cell.Value = $u.D
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
